# Auto-generated edit script: refresh the charging-station "not charged" report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write the text columns (A = station name, B = terminal name/id) as
# formulas that evaluate to string literals, so that numeric-looking terminal
# ids (e.g. "9176699400500301") are not silently reinterpreted as numbers.
$ws.Range("A2").Formula = '="长沙市开福区高岭香江国际城充电站建设项目"'
$ws.Range("B2").Formula = '="105号直流"'
$ws.Range("A3").Formula = '="长沙市开福区高岭香江国际城充电站建设项目"'
$ws.Range("B3").Formula = '="204号直流"'
$ws.Range("A4").Formula = '="长沙市开福区高岭香江国际城充电站建设项目"'
$ws.Range("B4").Formula = '="106号直流"'
$ws.Range("A5").Formula = '="长沙市开福区高岭香江国际城充电站建设项目"'
$ws.Range("B5").Formula = '="206号直流"'
$ws.Range("A6").Formula = '="长沙市开福区高岭香江国际城充电站建设项目"'
$ws.Range("B6").Formula = '="110号直流"'
$ws.Range("A7").Formula = '="长沙市开福区高岭香江国际城充电站建设项目"'
$ws.Range("B7").Formula = '="207号直流"'
$ws.Range("A8").Formula = '="长沙市开福区高岭香江国际城充电站建设项目"'
$ws.Range("B8").Formula = '="111号直流"'
$ws.Range("A9").Formula = '="长沙市开福区高岭香江国际城充电站建设项目"'
$ws.Range("B9").Formula = '="108号直流"'
$ws.Range("A10").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B10").Formula = '="9176699355900102"'
$ws.Range("A11").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B11").Formula = '="9176699400500403"'
$ws.Range("A12").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B12").Formula = '="9176699400501305"'
$ws.Range("A13").Formula = '="飞狐四方坪东区充电站"'
$ws.Range("B13").Formula = '="9176699442101001"'
$ws.Range("A14").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B14").Formula = '="9176699400500605"'
$ws.Range("A15").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B15").Formula = '="9176699400500204"'
$ws.Range("A16").Formula = '="飞狐四方坪东区充电站"'
$ws.Range("B16").Formula = '="9176699425700301"'
$ws.Range("A17").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B17").Formula = '="9176699400500301"'
$ws.Range("A18").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B18").Formula = '="9176699400500102"'
$ws.Range("A19").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B19").Formula = '="9176699400500604"'
$ws.Range("A20").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B20").Formula = '="9176699400500304"'
$ws.Range("A21").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B21").Formula = '="9176699400501205"'
$ws.Range("A22").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B22").Formula = '="9176699400500501"'
$ws.Range("A23").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B23").Formula = '="9176699400500502"'
$ws.Range("A24").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B24").Formula = '="9176699400500504"'
$ws.Range("A25").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B25").Formula = '="9176699400500303"'
$ws.Range("A26").Formula = '="飞狐四方坪东区充电站"'
$ws.Range("B26").Formula = '="9176699442100302"'
$ws.Range("A27").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B27").Formula = '="9176699400501303"'
$ws.Range("A28").Formula = '="飞狐四方坪西区充电站"'
$ws.Range("B28").Formula = '="9176699400500201"'

# Step 2: convert those formulas into plain static text values in-place
# (copy + paste-special values) so the saved cells keep their original cell
# style/number-format (General for A, several cells in B) while storing the
# content as text, exactly like the source data.
$textRange = $ws.Range("A2:B28")
$textRange.Copy() | Out-Null
$textRange.PasteSpecial(-4163) | Out-Null

# Step 3: write the date/time columns (C = last-seen-unplugged timestamp,
# D = report-generation timestamp) as plain numeric serial dates.
$ws.Range("C2").Value = 46043.552719907406
$ws.Range("D2").Value = 46044.328784722224
$ws.Range("C3").Value = 46043.554050925923
$ws.Range("D3").Value = 46044.328784722224
$ws.Range("C4").Value = 46043.624456018515
$ws.Range("D4").Value = 46044.328784722224
$ws.Range("C5").Value = 46043.649282407408
$ws.Range("D5").Value = 46044.328784722224
$ws.Range("C6").Value = 46043.666828703703
$ws.Range("D6").Value = 46044.328784722224
$ws.Range("C7").Value = 46043.706435185188
$ws.Range("D7").Value = 46044.328784722224
$ws.Range("C8").Value = 46043.70853009259
$ws.Range("D8").Value = 46044.328784722224
$ws.Range("C9").Value = 46043.7112037037
$ws.Range("D9").Value = 46044.328784722224
$ws.Range("C10").Value = 46043.112453703703
$ws.Range("D10").Value = 46044.310810185183
$ws.Range("C11").Value = 46043.377685185187
$ws.Range("D11").Value = 46044.310810185183
$ws.Range("C12").Value = 46043.565451388888
$ws.Range("D12").Value = 46044.310810185183
$ws.Range("C13").Value = 46043.573854166665
$ws.Range("D13").Value = 46044.310810185183
$ws.Range("C14").Value = 46043.576226851852
$ws.Range("D14").Value = 46044.310810185183
$ws.Range("C15").Value = 46043.576481481483
$ws.Range("D15").Value = 46044.310810185183
$ws.Range("C16").Value = 46043.576921296299
$ws.Range("D16").Value = 46044.310810185183
$ws.Range("C17").Value = 46043.579606481479
$ws.Range("D17").Value = 46044.310810185183
$ws.Range("C18").Value = 46043.5859375
$ws.Range("D18").Value = 46044.310810185183
$ws.Range("C19").Value = 46043.590740740743
$ws.Range("D19").Value = 46044.310810185183
$ws.Range("C20").Value = 46043.592974537038
$ws.Range("D20").Value = 46044.310810185183
$ws.Range("C21").Value = 46043.593657407408
$ws.Range("D21").Value = 46044.310810185183
$ws.Range("C22").Value = 46043.595810185187
$ws.Range("D22").Value = 46044.310810185183
$ws.Range("C23").Value = 46043.604733796295
$ws.Range("D23").Value = 46044.310810185183
$ws.Range("C24").Value = 46043.625358796293
$ws.Range("D24").Value = 46044.310810185183
$ws.Range("C25").Value = 46043.634027777778
$ws.Range("D25").Value = 46044.310810185183
$ws.Range("C26").Value = 46043.651747685188
$ws.Range("D26").Value = 46044.310810185183
$ws.Range("C27").Value = 46043.722638888888
$ws.Range("D27").Value = 46044.310810185183
$ws.Range("C28").Value = 46043.73709490741
$ws.Range("D28").Value = 46044.310810185183

# Step 4: restore the active selection to match the refreshed workbook.
$ws.Range("E23").Select() | Out-Null
